$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Wipe existing hyperlinks + cell content/formatting so we can rebuild the
#    sheet from scratch at the final target cell addresses (the insert/shift
#    operations in this runtime do not relocate hyperlinks correctly, so we
#    avoid relying on them).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()
$ws.Cells.Clear()

# ---------------------------------------------------------------------------
# 2. Values
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "method"
$ws.Range("C1").Value = "url"
$ws.Range("D1").Value = "headers"
$ws.Range("E1").Value = "params"
$ws.Range("F1").Value = "form_data"
$ws.Range("G1").Value = "data"
$ws.Range("H1").Value = "json"
$ws.Range("I1").Value = "expected_status_code"
$ws.Range("J1").Value = "dependencies"
$ws.Range("K1").Value = "dependency_extracts"

$ws.Range("A2").Value = "登录"
$ws.Range("B2").Value = "post"
$ws.Range("C2").Value = "https://api.hongrenju96.com/user/open/user/login"
$ws.Range("D2").Value = "Access-Domain=zt"
$ws.Range("E2").Value = ""
$ws.Range("G2").Value = "{'userName': '13288837480',
'password': '13288837480',
'userType': '-1',
'logPort': '0',
'logType': '1',
'domain': 'zt'}"
$ws.Range("I2").Value = 1001

$ws.Range("A3").Value = "更改密码"
$ws.Range("B3").Value = "post"
$ws.Range("C3").Value = "https://api.hongrenju96.com/api/admin/user/updateUserPassword"
$ws.Range("D3").Value = "Access-Domain=zt"
$ws.Range("E3").Value = ""
$ws.Range("G3").Value = "{'id':'c39bf60707f64c108bf1ed7ae8641e9f','newPassword':'13288837481','token':`${token}
}"
$ws.Range("I3").Value = 1001
$ws.Range("J3").Value = "登录"
$ws.Range("K3").Value = "登录:data.token"

$ws.Range("A4").Value = "知识库列表"
$ws.Range("B4").Value = "get"
$ws.Range("C4").Value = "https://api.hongrenju96.com/ymkEnterprise/pc/admin/enterpriseKnowledge/page?page=1&size=20"
$ws.Range("D4").Value = "Access-Domain=zt"
$ws.Range("E4").Value = ""
$ws.Range("I4").Value = 1001

$ws.Range("A5").Value = "新建知识库"
$ws.Range("B5").Value = "post"
$ws.Range("C5").Value = "https://api.hongrenju96.com/ymkEnterprise/pc/admin/enterpriseKnowledge/saveProduct"
$ws.Range("D5").Value = "Access-Domain=zt"
$ws.Range("G5").Value = '{"name":"袜子","description":"阿萨德","targetUsers":"20-30","keySellingPoints":"这是一双来自深渊的袜子","price":20,"id":null,"images":[],"token":"token"}'
$ws.Range("I5").Value = 1001

$ws.Range("A6").Value = "删除知识库"
$ws.Range("B6").Value = "post"
$ws.Range("C6").Value = "https://api.hongrenju96.com/ymkEnterprise/pc/admin/enterpriseKnowledge/delete"
$ws.Range("D6").Value = "Access-Domain=zt"

Write-Output "values set"
